# Daily attendance processing - 2025-11-28 07:02:55
# Updates the "Recorded By" (column G) values for specific rows in the
# "Session Analysis Results" sheet, rotating the comma-separated list of
# recorders left by one position (first recorder moved to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of row number -> new "Recorded By" value (column G)
$updates = [ordered]@{
    2   = "system, backup@backdoor.com, System"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    7   = "admin@admin.com, System"
    8   = "backup@backdoor.com, System"
    11  = "System, dnasr281@gmail.com"
    17  = "System, dnasr281@gmail.com"
    28  = "system, backup@backdoor.com, System"
    30  = "backup@backdoor.com, System"
    31  = "backup@backdoor.com, System"
    33  = "admin@admin.com, System"
    34  = "backup@backdoor.com, System"
    37  = "System, dnasr281@gmail.com"
    43  = "System, dnasr281@gmail.com"
    54  = "system, backup@backdoor.com, System"
    56  = "backup@backdoor.com, System"
    57  = "backup@backdoor.com, System"
    59  = "admin@admin.com, System"
    60  = "backup@backdoor.com, System"
    63  = "System, dnasr281@gmail.com"
    69  = "System, dnasr281@gmail.com"
    80  = "backup@backdoor.com, System"
    81  = "backup@backdoor.com, System"
    82  = "backup@backdoor.com, System"
    87  = "admin@admin.com, dnasr281@gmail.com"
    93  = "System, dnasr281@gmail.com"
    94  = "System, dnasr281@gmail.com"
    96  = "System, dnasr281@gmail.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    113 = "admin@admin.com, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    120 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    139 = "admin@admin.com, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    146 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
